$d = $word.ActiveDocument
$rng = $d.Content
$f = $rng.Find
$f.Text = "Pitch"
$f.Execute() | Out-Null
Write-Output "found text=[$($rng.Text)] start=$($rng.Start) end=$($rng.End)"
$xml = "<w:r xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:rPr><w:rFonts w:eastAsia='Times New Roman' w:cstheme='minorHAnsi'/><w:lang w:eastAsia='de-DE'/></w:rPr><w:lastRenderedPageBreak/><w:t>Pitch</w:t></w:r>"
$rng.InsertXML($xml)
